$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the German keyword translations for "child labor" and "development aid"
$ws.Range("D8").Value = "Kinderarbeit"
$ws.Range("D9").Value = "Entwicklungshilfe"

# Update the active selection to D14
$ws.Range("D14").Select()
